# edit.ps1
# Applies the "Doing Updates for Financials" commit:
#  - Inserts two new columns (D, E) before the existing quarterly data to
#    hold the two newest reporting quarters (period-end 2018-12-31 / 43465
#    and 2018-09-30 / 43373, encoded as Excel date serials).
#  - The pre-existing quarterly columns shift right from D:K to F:M.
#  - A handful of previously-reported quarterly figures were restated
#    (rows 42, 48, 49, 58, 91) - those corrected values are re-applied
#    to their (now shifted) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Insert two blank columns at D:E - this pushes the existing D:K data to F:M.
$ws.Columns("D:E").Insert()

# 2) Clone formatting (number format, font, blank/NA placeholders, etc.) for
#    the two new columns from column F (which now holds the data that used
#    to live in column D), so the new cells look like the rest of the table.
$ws.Range("F5:F102").Copy($ws.Range("D5:D102"))
$ws.Range("F5:F102").Copy($ws.Range("E5:E102"))

# 3) Populate the two new columns with the newest-quarter figures.
#    Each entry is (row, D-value, E-value); $null means "leave the cell blank"
#    (header / section-separator rows).
$newQuarterData = @(
    @(7, 43465, 43373),
    @(8, 2381000, 1512000),
    @(9, 832000, 513000),
    @(10, 1549000, 999000),
    @(11, $null, $null),
    @(12, 325000, 263000),
    @(13, 0, 0),
    @(14, 0, 40000),
    @(15, 0, 0),
    @(16, $null, $null),
    @(17, 1687000, 1287000),
    @(18, 694000, 225000),
    @(19, $null, $null),
    @(20, 18000, 20000),
    @(21, 1087000, 395000),
    @(22, 22000, 33000),
    @(23, 690000, 212000),
    @(24, 324000, -48000),
    @(25, 0, 0),
    @(26, 366000, 260000),
    @(27, 366000, 260000),
    @(28, 0, 0),
    @(29, 285000, "NA"),
    @(30, 0, 0),
    @(31, 0, 0),
    @(32, -18000, -20000),
    @(33, 651000, 260000),
    @(34, 0, 0),
    @(35, 651000, 260000),
    @(38, 43465, 43373),
    @(39, $null, $null),
    @(40, $null, $null),
    @(41, 4225000, 3308000),
    @(42, 0, 0),
    @(43, 1035000, 641000),
    @(44, 43000, 174000),
    @(45, 803000, 849000),
    @(46, 6106000, 4972000),
    @(47, 0, 0),
    @(48, 282000, 281000),
    @(49, 10562000, 10763000),
    @(50, 0, 0),
    @(51, 0, 0),
    @(52, 885000, 697000),
    @(53, 0, 0),
    @(54, 17835000, 16713000),
    @(55, $null, $null),
    @(56, $null, $null),
    @(57, 253000, 312000),
    @(58, 0, 0),
    @(59, 2389000, 2070000),
    @(60, 2642000, 2382000),
    @(61, 2671000, 2670000),
    @(62, 1165000, 1002000),
    @(63, 0, 0),
    @(64, 0, 0),
    @(65, 0, 0),
    @(66, 6478000, 6054000),
    @(67, $null, $null),
    @(68, 0, 0),
    @(69, 0, 0),
    @(70, 0, 0),
    @(71, 0, 0),
    @(72, 6558000, 5907000),
    @(73, 0, 0),
    @(74, 0, 0),
    @(75, 0, 0),
    @(76, 11357000, 10659000),
    @(77, 0, 0),
    @(80, 43465, 43373),
    @(81, 651000, 260000),
    @(82, $null, $null),
    @(83, 375000, 150000),
    @(84, 0, 0),
    @(85, 0, 0),
    @(86, 0, 0),
    @(87, 0, 0),
    @(88, 0, 0),
    @(89, 999000, 253000),
    @(90, $null, $null),
    @(91, -34000, -36000),
    @(92, 0, 0),
    @(93, 0, 0),
    @(94, -70000, -36000),
    @(95, $null, $null),
    @(96, 0, 0),
    @(97, 0, 0),
    @(98, 0, 0),
    @(99, 0, 0),
    @(100, 0, -1770000),
    @(101, -16000, 4000),
    @(102, 913000, -1549000)
)

foreach ($entry in $newQuarterData) {
    $r = $entry[0]
    $dVal = $entry[1]
    $eVal = $entry[2]
    if ($null -ne $dVal) {
        $ws.Cells.Item($r, 4).Value = $dVal
    } else {
        $ws.Cells.Item($r, 4).Value = ""
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($r, 5).Value = $eVal
    } else {
        $ws.Cells.Item($r, 5).Value = ""
    }
}

# 4) Re-apply the handful of restated historical values (now shifted into
#    columns F:J) that changed as part of this data refresh.
$restatements = @(
    @(42, "F", 0),
    @(42, "G", 0),
    @(42, "H", 0),
    @(42, "I", 0),
    @(42, "J", 0),
    @(48, "F", 281000),
    @(48, "G", 286000),
    @(48, "H", 294000),
    @(48, "I", 254000),
    @(48, "J", 246000),
    @(49, "F", 10804000),
    @(49, "G", 10843000),
    @(49, "H", 10955000),
    @(49, "I", 11170000),
    @(49, "J", 11346000),
    @(58, "F", 0),
    @(58, "G", 0),
    @(58, "H", 0),
    @(58, "I", 0),
    @(58, "J", 0),
    @(91, "I", -34000),
    @(91, "J", -31000)
)

foreach ($entry in $restatements) {
    $r = $entry[0]
    $col = $entry[1]
    $val = $entry[2]
    $ws.Range("$col$r").Value = $val
}
